$wb = $excel.ActiveWorkbook

# --- Commodity sheet: remove WindOff and WindOn commodity rows (rows 2 and 3) ---
$wsCommodity = $wb.Worksheets.Item("Commodity")
$wsCommodity.Rows.Item(2).Resize(2).Delete() | Out-Null

# --- Process sheet: remove Wind (onshore) and Wind (offshore) process rows (rows 2 and 3) ---
$wsProcess = $wb.Worksheets.Item("Process")
$wsProcess.Rows.Item(2).Resize(2).Delete() | Out-Null

# --- Process-Commodity sheet: remove the 4 rows describing Wind (onshore)/Wind (offshore)
#     In/Out commodity ratios (rows 2-5) ---
$wsProcessCommodity = $wb.Worksheets.Item("Process-Commodity")
$wsProcessCommodity.Rows.Item(2).Resize(4).Delete() | Out-Null

# --- SupIm sheet: remove the EU27.WindOff and EU27.WindOn columns (columns B and C) ---
$wsSupIm = $wb.Worksheets.Item("SupIm")
$wsSupIm.Columns.Item(2).Resize(1, 2).Delete() | Out-Null

# --- Make SupIm the active sheet/selection, matching the saved view state ---
$wsSupIm.Select()
$wsSupIm.Range("B1:C1048576").Select()
